$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("D2").Value = "'28.782.83"
$ws.Range("E2").Value = "  -1.80%  "

# Row 3
$ws.Range("D3").Value = "'1.827.88"
$ws.Range("E3").Value = "  -1.75%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'244.00"
$ws.Range("E5").Value = "  +0.72%  "

# Row 6
$ws.Range("D6").Value = "'0.6873"
$ws.Range("E6").Value = "  -1.82%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.07614"
$ws.Range("E8").Value = "  -3.25%  "

# Row 9
$ws.Range("D9").Value = "'0.3030"
$ws.Range("E9").Value = "  -3.00%  "

# Row 10
$ws.Range("D10").Value = "'23.09"
$ws.Range("E10").Value = "  -4.33%  "

# Row 11
$ws.Range("D11").Value = "'0.07784"
$ws.Range("E11").Value = "  -0.20%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.828.61"
$ws.Range("E12").Value = "  -1.91%  "

# Row 13
$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'92.47"
$ws.Range("E13").Value = "  +0.12%  "

# Row 14
$ws.Range("D14").Value = "'5.070"
$ws.Range("E14").Value = "  -1.36%  "

# Row 15
$ws.Range("D15").Value = "'0.6749"
$ws.Range("E15").Value = "  -2.78%  "

# Row 16
$ws.Range("D16").Value = "'6.429"
$ws.Range("E16").Value = "  -1.03%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'28.818.56"
$ws.Range("E17").Value = "  -1.75%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000008199"
$ws.Range("E18").Value = "  -3.73%  "

# Row 19
$ws.Range("D19").Value = "'240.67"
$ws.Range("E19").Value = "  -3.12%  "

# Row 20
$ws.Range("D20").Value = "'2.073.02"
$ws.Range("E20").Value = "  -1.99%  "

# Row 21
$ws.Range("D21").Value = "'12.63"
$ws.Range("E21").Value = "  -2.59%  "

# Row 22
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$ws.Range("D23").Value = "'7.405"
$ws.Range("E23").Value = "  -1.90%  "

# Row 24
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("D25").Value = "'0.1476"
$ws.Range("E25").Value = "  -3.86%  "

# Row 26
$ws.Range("D26").Value = "'162.24"
$ws.Range("E26").Value = "  +0.94%  "

# Row 27
$ws.Range("D27").Value = "'8.706"
$ws.Range("E27").Value = "  -2.53%  "

# Row 28
$ws.Range("E28").Value = "  -2.76%  "

# Row 29
$ws.Range("E29").Value = "  -2.44%  "

# Row 30
$ws.Range("D30").Value = "'4.203"
$ws.Range("E30").Value = "  -1.80%  "

# Row 31
$ws.Range("D31").Value = "'4.143"
$ws.Range("E31").Value = "  -2.53%  "

# Row 32
$ws.Range("D32").Value = "'1.185"
$ws.Range("E32").Value = "  -1.78%  "

# Row 33
$ws.Range("D33").Value = "'0.05083"
$ws.Range("E33").Value = "  -2.88%  "

# Row 34
$ws.Range("D34").Value = "'0.7670"
$ws.Range("E34").Value = "  +2.28%  "

# Row 35
$ws.Range("D35").Value = "'1.836"
$ws.Range("E35").Value = "  -2.47%  "

# Row 36
$ws.Range("D36").Value = "'1.131"
$ws.Range("E36").Value = "  -3.75%  "

# Row 37
$ws.Range("D37").Value = "'2.693"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$ws.Range("D38").Value = "'0.01853"
$ws.Range("E38").Value = "  -0.54%  "

# Row 39
$ws.Range("D39").Value = "'1.242.47"
$ws.Range("E39").Value = "  -1.84%  "

# Row 40
$ws.Range("D40").Value = "'2.697"
$ws.Range("E40").Value = "  -1.69%  "

# Row 41
$ws.Range("D41").Value = "'0.9537"
$ws.Range("E41").Value = "  +6.18%  "

# Row 42
$ws.Range("D42").Value = "'5.954"
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").Value = "'106.60"
$ws.Range("E43").Value = "  -3.85%  "

# Row 44
$ws.Range("E44").Value = "  +0.04%  "

# Row 45
$ws.Range("D45").Value = "'9.647"
$ws.Range("E45").Value = "  +0.91%  "

# Row 46
$ws.Range("E46").Value = "  -1.52%  "

# Row 47
$ws.Range("D47").Value = "'1.974.63"
$ws.Range("E47").Value = "  -2.05%  "

# Row 48
$ws.Range("D48").Value = "'0.5152"
$ws.Range("E48").Value = "  -0.62%  "

# Row 49
$ws.Range("D49").Value = "'63.37"
$ws.Range("E49").Value = "  -8.90%  "

# Row 50
$ws.Range("D50").Value = "'1.734"
$ws.Range("E50").Value = "  -2.50%  "

# Row 51
$ws.Range("D51").Value = "'6.910"
$ws.Range("E51").Value = "  -1.21%  "
